$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Copy formatting (style) from the last existing data row (577) to the new rows (578-600)
$ws.Range("A577:C577").Copy()
$ws.Range("A578:C600").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(578, 1).Value = "cs"
$ws.Cells.Item(578, 2).Value = "lab.setup.preview"
$ws.Cells.Item(578, 3).Value = "Náhled setupu"

$ws.Cells.Item(579, 1).Value = "cs"
$ws.Cells.Item(579, 2).Value = "lab.setup.deleted.success"
$ws.Cells.Item(579, 3).Value = "Setup byl úspěšně odstraněn."

$ws.Cells.Item(580, 1).Value = "cs"
$ws.Cells.Item(580, 2).Value = "lab.setup.button.index"
$ws.Cells.Item(580, 3).Value = "Detail setupu"

$ws.Cells.Item(581, 1).Value = "cs"
$ws.Cells.Item(581, 2).Value = "lab.setup.button.edit"
$ws.Cells.Item(581, 3).Value = "Editovat"

$ws.Cells.Item(582, 1).Value = "cs"
$ws.Cells.Item(582, 2).Value = "lab.setup.button.delete"
$ws.Cells.Item(582, 3).Value = "Odstranit setup"

$ws.Cells.Item(583, 1).Value = "cs"
$ws.Cells.Item(583, 2).Value = "lab.setup.button.delete.confirm.title"
$ws.Cells.Item(583, 3).Value = "Odstranit setup"

$ws.Cells.Item(584, 1).Value = "cs"
$ws.Cells.Item(584, 2).Value = "lab.setup.button.delete.confirm"
$ws.Cells.Item(584, 3).Value = "Opravdu si přejete odstranit vybraný setup? Tímto odstraníte veškerá data s ním spojená, můžete tak změnit i statistiky vapování a další vedlejší efekty. Použijte pouze pokud jste si naprosto jisti, co děláte; tuto akci nelze vzít zpět."
$ws.Rows.Item(584).RowHeight = 39

$ws.Cells.Item(585, 1).Value = "cs"
$ws.Cells.Item(585, 2).Value = "lab.setup.button.delete.confirm.ok"
$ws.Cells.Item(585, 3).Value = "Odstranit setup"

$ws.Cells.Item(586, 1).Value = "cs"
$ws.Cells.Item(586, 2).Value = "lab.setup.preview.preview.title"
$ws.Cells.Item(586, 3).Value = "Detail setupu"

$ws.Cells.Item(587, 1).Value = "cs"
$ws.Cells.Item(587, 2).Value = "lab.setup.preview.preview.subtitle"
$ws.Cells.Item(587, 3).Value = "Zde naleznete veškeré informace o vybraném setupu."

$ws.Cells.Item(588, 1).Value = "cs"
$ws.Cells.Item(588, 2).Value = "lab.setup.preview.name"
$ws.Cells.Item(588, 3).Value = "Jméno"

$ws.Cells.Item(589, 1).Value = "cs"
$ws.Cells.Item(589, 2).Value = "lab.setup.preview.description"
$ws.Cells.Item(589, 3).Value = "Popis"

$ws.Cells.Item(590, 1).Value = "cs"
$ws.Cells.Item(590, 2).Value = "lab.setup.preview.atomizer"
$ws.Cells.Item(590, 3).Value = "Atomizér"

$ws.Cells.Item(591, 1).Value = "cs"
$ws.Cells.Item(591, 2).Value = "lab.setup.preview.mod"
$ws.Cells.Item(591, 3).Value = "Mod"

$ws.Cells.Item(592, 1).Value = "cs"
$ws.Cells.Item(592, 2).Value = "lab.setup.preview.coil"
$ws.Cells.Item(592, 3).Value = "Spirálka"

$ws.Cells.Item(593, 1).Value = "cs"
$ws.Cells.Item(593, 2).Value = "lab.setup.preview.cotton"
$ws.Cells.Item(593, 3).Value = "Vata"

$ws.Cells.Item(594, 1).Value = "cs"
$ws.Cells.Item(594, 2).Value = "lab.setup.preview.ohm"
$ws.Cells.Item(594, 3).Value = "Odpor"

$ws.Cells.Item(595, 1).Value = "cs"
$ws.Cells.Item(595, 2).Value = "lab.setup.index.title"
$ws.Cells.Item(595, 3).Value = "Detail setupu"

$ws.Cells.Item(596, 1).Value = "cs"
$ws.Cells.Item(596, 2).Value = "lab.setup.index.preview.title"
$ws.Cells.Item(596, 3).Value = "Detail setupu"

$ws.Cells.Item(597, 1).Value = "cs"
$ws.Cells.Item(597, 2).Value = "lab.setup.index.preview.subtitle"
$ws.Cells.Item(597, 3).Value = "Zde naleznete veškeré informace o vybraném setupu."

$ws.Cells.Item(598, 1).Value = "cs"
$ws.Cells.Item(598, 2).Value = "lab.setup.edit.title"
$ws.Cells.Item(598, 3).Value = "Editace setupu"

$ws.Cells.Item(599, 1).Value = "cs"
$ws.Cells.Item(599, 2).Value = "lab.setup.edit.subtitle"
$ws.Cells.Item(599, 3).Value = "Pokud je třeba něco poladit…"

$ws.Cells.Item(600, 1).Value = "cs"
$ws.Cells.Item(600, 2).Value = "lab.setup.link.button"
$ws.Cells.Item(600, 3).Value = "Detail setupu"

# Match the final selection state from the authored edit
$ws.Range("C585").Select()
